# Reorders the course-requirement list in the "Requisitos" section of the
# document (LOB1223.docx). The set of 25 requirement lines is unchanged;
# only their order within the single ListBullet paragraph changes.

$d = $word.ActiveDocument

# Locate the "Requisitos" heading paragraph, then the list paragraph that
# immediately follows it (rather than hard-coding a paragraph index).
$count = $d.Paragraphs.Count
$listParaIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $heading = $d.Paragraphs.Item($i).Range.Text.Trim()
    if ($heading -eq "Requisitos") {
        $listParaIndex = $i + 1
        break
    }
}

if ($listParaIndex -eq -1) {
    throw "Could not find the 'Requisitos' heading paragraph"
}

$listPara = $d.Paragraphs.Item($listParaIndex)

# Replace the paragraph's runs (excluding the trailing paragraph mark) with
# the reordered set of runs via a flat-OPC OOXML fragment. This keeps each
# course requirement as its own <w:r> (with a <w:br/> between entries),
# matching the original document's run layout, and preserves the
# paragraph's own formatting (ListBullet style) since we don't touch the
# paragraph mark itself.
$targetRange = $d.Range($listPara.Range.Start, $listPara.Range.End - 1)

$flatOpc = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>LOB1053 -  Física III  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1004 -  Cálculo II  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1006 -  Cálculo IV  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1009 -  Leitura e Interpretação de Desenho Técnico  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1003 -  Cálculo I  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1040 -  Laboratório de Eletricidade  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1041 -  Física Experimental II  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1039 -  Física Experimental III  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOQ4097 -  Fundamentos de Química para Engenharia I (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1024 -  Mecânica  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1012 -  Estatística  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1268 -  Leitura, Escrita e Comunicação Científica  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOM3081 -  Introdução à Mecânica dos Sólidos  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1036 -  Geometria Analítica  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1018 -  Física I  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1270 -  Química Experimental Aplicada  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1019 -  Física II  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1011 -  Eletricidade Aplicada  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1037 -  Álgebra Linear  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1038 -  Física Experimental I  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1052 -  Cálculo III  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOQ4233 -  Gestão de Negócios  (Requisito fraco)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$targetRange.InsertXML($flatOpc)
